$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45 / 46: EnergySwap <-> Decentraland swap with updated price/volume
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6023"
$ws.Range("E45").Value = "'  +0.70%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'13.19"
$ws.Range("E46").Value = "'  -0.79%  "

# Price / volume updates
$ws.Range("D2").Value = "'28.062.52"
$ws.Range("E2").Value = "'  +0.03%  "
$ws.Range("D3").Value = "'1.873.74"
$ws.Range("E3").Value = "'  -1.67%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "'  -0.31%  "
$ws.Range("D5").Value = "'312.80"
$ws.Range("E5").Value = "'  +0.02%  "
$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "'  -0.32%  "
$ws.Range("D7").Value = "'0.5045"
$ws.Range("E7").Value = "'  -0.08%  "
$ws.Range("D8").Value = "'0.3840"
$ws.Range("E8").Value = "'  -2.05%  "
$ws.Range("D9").Value = "'0.09008"
$ws.Range("E9").Value = "'  -5.99%  "
$ws.Range("D10").Value = "'1.120"
$ws.Range("E10").Value = "'  -1.23%  "
$ws.Range("D11").Value = "'41.64"
$ws.Range("E11").Value = "'  -0.91%  "
$ws.Range("D12").Value = "'6.377"
$ws.Range("E12").Value = "'  +0.01%  "
$ws.Range("D13").Value = "'20.75"
$ws.Range("E13").Value = "'  -0.07%  "
$ws.Range("D14").Value = "'1.872.07"
$ws.Range("E14").Value = "'  -1.00%  "
$ws.Range("D15").Value = "'7.258"
$ws.Range("E15").Value = "'  -0.54%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = "'  -0.23%  "
$ws.Range("D17").Value = "'0.00001107"
$ws.Range("E17").Value = "'  -0.65%  "
$ws.Range("D18").Value = "'91.29"
$ws.Range("E18").Value = "'  -0.88%  "
$ws.Range("D19").Value = "'0.06644"
$ws.Range("E19").Value = "'  +0.71%  "
$ws.Range("D20").Value = "'18.22"
$ws.Range("E20").Value = "'  +2.20%  "
$ws.Range("D21").Value = "'0.9993"
$ws.Range("E21").Value = "'  -0.23%  "
$ws.Range("D22").Value = "'6.132"
$ws.Range("E22").Value = "'  -1.19%  "
$ws.Range("D23").Value = "'28.068.35"
$ws.Range("E23").Value = "'  -0.13%  "
$ws.Range("E24").Value = "'  +2.26%  "
$ws.Range("D25").Value = "'2.257"
$ws.Range("E25").Value = "'  -1.98%  "
$ws.Range("D26").Value = "'2.081.97"
$ws.Range("E26").Value = "'  -1.72%  "
$ws.Range("D27").Value = "'2.534"
$ws.Range("E27").Value = "'  -4.84%  "
$ws.Range("D28").Value = "'20.80"
$ws.Range("E28").Value = "'  +0.06%  "
$ws.Range("D29").Value = "'157.09"
$ws.Range("E29").Value = "'  -0.27%  "
$ws.Range("D30").Value = "'127.37"
$ws.Range("E30").Value = "'  +0.35%  "
$ws.Range("E31").Value = "'  +0.40%  "
$ws.Range("D32").Value = "'1.061"
$ws.Range("E32").Value = "'  -2.03%  "
$ws.Range("D33").Value = "'5.610"
$ws.Range("E33").Value = "'  +0.02%  "
$ws.Range("D34").Value = "'3.597"
$ws.Range("E34").Value = "'  -0.53%  "
$ws.Range("D35").Value = "'9.444"
$ws.Range("E35").Value = "'  -1.85%  "
$ws.Range("D36").Value = "'0.06584"
$ws.Range("E36").Value = "'  -0.21%  "
$ws.Range("E37").Value = "'  -0.91%  "
$ws.Range("E38").Value = "'  +0.59%  "
$ws.Range("D39").Value = "'1.290"
$ws.Range("E39").Value = "'  +1.06%  "
$ws.Range("D40").Value = "'1.210"
$ws.Range("E40").Value = "'  -1.67%  "
$ws.Range("D41").Value = "'0.6392"
$ws.Range("E41").Value = "'  +0.94%  "
$ws.Range("D42").Value = "'11.50"
$ws.Range("E42").Value = "'  +1.29%  "
$ws.Range("D43").Value = "'4.921"
$ws.Range("E43").Value = "'  -1.15%  "
$ws.Range("D44").Value = "'0.9990"
$ws.Range("E44").Value = "'  -0.27%  "
$ws.Range("E47").Value = "'  -0.22%  "
$ws.Range("D48").Value = "'3.664"
$ws.Range("E48").Value = "'  -1.62%  "
$ws.Range("D49").Value = "'1.237"
$ws.Range("E49").Value = "'  +4.71%  "
$ws.Range("D50").Value = "'1.999"
$ws.Range("E50").Value = "'  -1.12%  "
$ws.Range("D51").Value = "'121.03"
$ws.Range("E51").Value = "'  -1.34%  "
